# "finalisation travail avant vacances hiver"
# Add a new row (row 10) to Sheet1, re-using the same layout / red-font
# style as the other "alpha fixé ..." rows (A5:A8), and move the active
# selection to D10 to match the saved workbook state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New label in column A, styled like A5 ("alpha fixé et  gamma à optimiser")
$ws.Range("A10").Value = "alpha fixé"
$ws.Range("A10").Font.Color = $ws.Range("A5").Font.Color

# Numeric inputs
$ws.Range("B10").Value = 0.18
$ws.Range("E10").Value = 0.258

# Formulas (mirrors F5/G5 but without the extra parentheses around the constant)
$ws.Range("F10").Formula = "=ABS(E10-0.378)"
$ws.Range("G10").Formula = "=ABS(E10-0.138)"

$ws.Range("K10").Value = 1524.4584535804599
$ws.Range("L10").Value = 1
$ws.Range("M10").Formula = "=K10+2*L10"

# Match the saved selection state
$ws.Range("D10").Select()
